# Delete row 109 (AkNr "GG" / "Grundgesetz für die Bundesrepublik Deutschland" /
# "Basic Law for the Federal Republic of Germany") from the abbreviations list.
# Deleting the entire row shifts all subsequent rows up by one, which matches
# the target diff (old row N+1 becomes new row N for 110..190, and the sheet
# shrinks from A1:C190 to A1:C189).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dic_Abkürzungen")

$ws.Rows.Item(109).Delete()
